{"js": "// Title: \"...\u4e00\u5916\u4e00\u5167\u6bd42\" -> \"...\u8a71\u8a9e\u9738\u6b0a\"\n// Body: collapse everything after the title into a single question\n// paragraph reading \"(1) 21\u4e16\u7d00\u8cc7\u672c\u8ad6\" (images + questions (2)-(4) removed).\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst items = paras.items;\n\n// 1) Update the title paragraph's text.\nitems[0].insertText(\n  \"\u6578\u5b78 - \u61c9\u7528\u984c - \u5178\u578b\u61c9\u7528\u984c - \u71d5\u5c3e\u5b9a\u7406 - \u8a71\u8a9e\u9738\u6b0a\",\n  \"Replace\"\n);\n\n// 2) Update the first question paragraph's text (this also drops the\n//    trailing <w:br/> that used to live in that run, since the whole\n//    paragraph's text is replaced).\nitems[1].insertText(\"(1) 21\u4e16\u7d00\u8cc7\u672c\u8ad6\", \"Replace\");\n\n// 3) Remove every remaining paragraph (the image paragraph, (2), (3),\n//    the two-image paragraph, and (4)).\nfor (let i = items.length - 1; i >= 2; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Title: \"...\u4e00\u5916\u4e00\u5167\u6bd42\" -> \"...\u8a71\u8a9e\u9738\u6b0a\"\n# Body: collapse everything after the title into a single question\n# paragraph reading \"(1) 21\u4e16\u7d00\u8cc7\u672c\u8ad6\" (images + questions (2)-(4) removed).\n\n$d = $word.ActiveDocument\n\n# Remove every paragraph after the first question paragraph (the image\n# paragraph, (2), (3), the two-image paragraph, and (4)), working from the\n# end so earlier indices stay valid.\nfor ($i = $d.Paragraphs.Count; $i -ge 3; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Update the title paragraph's text.\n$d.Paragraphs.Item(1).Range.Text = \"\u6578\u5b78 - \u61c9\u7528\u984c - \u5178\u578b\u61c9\u7528\u984c - \u71d5\u5c3e\u5b9a\u7406 - \u8a71\u8a9e\u9738\u6b0a\"\n\n# Update the first question paragraph's text (replacing the whole range\n# also drops the trailing manual line break that used to live there).\n$d.Paragraphs.Item(2).Range.Text = \"(1) 21\u4e16\u7d00\u8cc7\u672c\u8ad6\"\n"}
